$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Delete the empty placeholder "Anleitung" sheet (last sheet, sheetId 12)
$emptySheet = $wb.Worksheets.Item("Anleitung")
$emptySheet.Delete()

# Rename "Anleitung1" to "Anleitung" now that the name is free
$ws = $wb.Worksheets.Item("Anleitung1")
$ws.Name = "Anleitung"
